$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 300
$ws1.Range("F5").Value = 5336
$ws1.Range("F6").Value = 613
$ws1.Range("F7").Value = 451
$ws1.Range("F8").Value = 239
$ws1.Range("F9").Value = 1098
$ws1.Range("F10").Value = 289
$ws1.Range("F11").Value = 155
$ws1.Range("F13").Value = 746
$ws1.Range("F17").Value = 188
$ws1.Range("F19").Value = 373
$ws1.Range("F20").Value = 6150
$ws1.Range("G20").Value = 78
$ws1.Range("F21").Value = 49
$ws1.Range("F22").Value = 47
$ws1.Range("F24").Value = 7101
$ws1.Range("F27").Value = 3269
$ws1.Range("F28").Value = 384
$ws1.Range("F29").Value = 772
$ws1.Range("F30").Value = 4467
$ws1.Range("F32").Value = 144
$ws1.Range("F34").Value = 1196
$ws1.Range("F35").Value = 109
$ws1.Range("F36").Value = 32
$ws1.Range("F38").Value = 952
$ws1.Range("F39").Value = 1200
$ws1.Range("F40").Value = 2064

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 232

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 232
$ws4.Range("F7").Value = 300
$ws4.Range("F8").Value = 5336
$ws4.Range("F9").Value = 613
$ws4.Range("F10").Value = 451
$ws4.Range("F11").Value = 239
$ws4.Range("F12").Value = 1098
$ws4.Range("F13").Value = 289
$ws4.Range("F14").Value = 155
$ws4.Range("F16").Value = 746
$ws4.Range("F21").Value = 188
$ws4.Range("F23").Value = 373
$ws4.Range("F24").Value = 6150
$ws4.Range("G24").Value = 78
$ws4.Range("F25").Value = 49
$ws4.Range("F26").Value = 47
$ws4.Range("F28").Value = 7101
$ws4.Range("F31").Value = 3269
$ws4.Range("F32").Value = 384
$ws4.Range("F33").Value = 772
$ws4.Range("F34").Value = 4467
$ws4.Range("F37").Value = 144
$ws4.Range("F39").Value = 1196
$ws4.Range("F40").Value = 109
$ws4.Range("F41").Value = 32
$ws4.Range("F43").Value = 952
$ws4.Range("F44").Value = 1200
$ws4.Range("F46").Value = 2064
